$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A entirely; this shifts columns B:F left to A:E,
# matching the diff (old column A values duplicated old column F and
# were removed, everything else shifted one column to the left).
$ws.Range("A1").EntireColumn.Delete()
